# Add three new wishlist items (rows 19-21) to the "Valentin" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Name) for the three new rows, filled first -------------------
$ws.Range("A19").Value = 'GN Red & Black ''HUD'' Mouse Mat'
$ws.Range("A20").Value = 'GN ''Charge'' Red/Black Microfiber Cloth Mousepad'
$ws.Range("A21").Value = 'GN Drink ''Debug'' Coaster Pack'

# --- Column C (Link), filled second -----------------------------------------
$ws.Range("C19").Value = 'https://store.gamersnexus.net/products/gn-red-black-hud-mouse-mat'
$ws.Range("C20").Value = 'https://store.gamersnexus.net/products/gn-charge-redblack-mousepad'
$ws.Range("C21").Value = 'https://store.gamersnexus.net/products/gn-drink-debug-coaster-pack-4-custom-3d-coasters-100x100mm-4x4'

# --- Column B (Image), filled third -----------------------------------------
$ws.Range("B19").Value = 'https://images.squarespace-cdn.com/content/v1/58c839976b8f5ba6ce6a8105/1639861410250-81GPR6I1BSEPA4D21TWM/top-down-plants-only2.jpg?format=750w'
$ws.Range("B20").Value = 'https://images.squarespace-cdn.com/content/v1/58c839976b8f5ba6ce6a8105/1621098515020-CZO4LA6HWNQRH4H9C5C9/charge-mousepad_desk-right.jpg?format=750w'
$ws.Range("B21").Value = 'https://images.squarespace-cdn.com/content/v1/58c839976b8f5ba6ce6a8105/1666119266873-0IJXCW0IUAFYQNK2W89O/coaster-pack-4-spaced2.jpg?format=750w'

# --- Column D (Price), filled last ------------------------------------------
# These look like plain currency numbers, so Excel would normally convert them
# to numeric values with a currency number format. Force the cells to Text
# first so the values round-trip as the literal strings "$29.99" / "$19.99",
# then restore the cell format so no extra formatting lingers on the cells.
$priceRange = $ws.Range("D19:D21")
$priceRange.NumberFormat = "@"

$ws.Range("D19").Value = '$29.99'
$ws.Range("D20").Value = '$19.99'
$ws.Range("D21").Value = '$29.99'

$priceRange.Style = "Normal"

# --- Update the active selection to match the final cursor position ---------
[void]$ws.Range("D20").Select()
